$wb = $excel.ActiveWorkbook

$wsMax = $wb.Worksheets.Item("max_capacity")
$wsReorder = $wb.Worksheets.Item("reorder_level")

# --- max_capacity sheet: forced re-order level raised from 100 to max capacity of 20 ---
$wsMax.Range("B2:C5").Value = 20

# --- reorder_level sheet: first spare part's reorder point lowered from 50/40 to 10/10 ---
$wsReorder.Range("B2:C2").Value = 10

# --- Update selection on reorder_level (no longer the active tab) ---
$wsReorder.Range("C2").Select()

# --- max_capacity becomes the active tab/selection ---
$wsMax.Activate()
$wsMax.Range("C5").Select()
